# UST to EST, ticker_list.xlsx to predict from list.
# Remove "elonmusk" from the user_names sheet's "controversial" column (C),
# shifting the remaining entries (JeffBezos, BillGates) up by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Manually shift column C values up starting at the row that held "elonmusk",
# instead of using Range.Delete (which in this environment shifts whole rows
# across all columns rather than just the targeted column).
$ws.Range("C2").Value = $ws.Range("C3").Value2
$ws.Range("C3").Value = $ws.Range("C4").Value2
$ws.Range("C4").ClearContents()

# Match the saved selection/active cell recorded in the workbook after the edit.
$ws.Range("C5").Select()
